$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 435, shifting existing rows 435:530 down by one.
$ws.Rows.Item(435).Insert()

# Populate the new row 435 with the new weekly price record.
$ws.Cells.Item(435, 1).Value = 10
$ws.Cells.Item(435, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(435, 3).Value = "La Araucanía"
$ws.Cells.Item(435, 4).Value = 45244
$ws.Cells.Item(435, 5).Value = 9
$ws.Cells.Item(435, 6).Value = "Fruta"
$ws.Cells.Item(435, 7).Value = 100102
$ws.Cells.Item(435, 8).Value = "Cítricos"
$ws.Cells.Item(435, 9).Value = 100102006
$ws.Cells.Item(435, 10).Value = "Pomelo"
$ws.Cells.Item(435, 11).Value = "Start Ruby"
$ws.Cells.Item(435, 12).Value = "Primera"
$ws.Cells.Item(435, 13).Value = 90
$ws.Cells.Item(435, 14).Value = 14000
$ws.Cells.Item(435, 15).Value = 16000
$ws.Cells.Item(435, 16).Value = 14889
$ws.Cells.Item(435, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(435, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(435, 19).Value = 993
$ws.Cells.Item(435, 20).Value = 15
